# "ajuste: corrigindo as categorias"
# - adds an "Idade ignorada" age-bracket column (S) and a "Total" column (T)
# - adds a new "Outros" disease-group row and a grand-"Total" row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns S (Idade ignorada) and T (Total) ---
$ws.Range("S1").Value = "Idade ignorada"
$ws.Range("T1").Value = "Total"

# --- Existing rows 2-6: fill in the new Total column (S stays blank) ---
$ws.Range("T2").Value = 2237
$ws.Range("T3").Value = 286
$ws.Range("T4").Value = 916
$ws.Range("T5").Value = 239
$ws.Range("T6").Value = 1486

# --- New row 7: "Outros" ---
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 138
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 39
$ws.Range("F7").Value = 65
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = 90
$ws.Range("I7").Value = 64
$ws.Range("J7").Value = 81
$ws.Range("K7").Value = 93
$ws.Range("L7").Value = 109
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 131
$ws.Range("O7").Value = 145
$ws.Range("P7").Value = 144
$ws.Range("Q7").Value = 186
$ws.Range("R7").Value = 617
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 2108

# --- New row 8: "Total" ---
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 156
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 46
$ws.Range("F8").Value = 74
$ws.Range("G8").Value = 96
$ws.Range("H8").Value = 121
$ws.Range("I8").Value = 129
$ws.Range("J8").Value = 163
$ws.Range("K8").Value = 232
$ws.Range("L8").Value = 326
$ws.Range("M8").Value = 439
$ws.Range("N8").Value = 574
$ws.Range("O8").Value = 701
$ws.Range("P8").Value = 776
$ws.Range("Q8").Value = 814
$ws.Range("R8").Value = 2609
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 7272
